# Unmerge the four merged "Thời gian" cell pairs in column E and fill the
# previously-blank half of each pair with the same value (and matching
# formatting) as the filled half, then move the active selection to G4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mergedPairs = @(
    @("E8", "E9"),
    @("E10", "E11"),
    @("E16", "E17"),
    @("E23", "E24")
)

foreach ($pair in $mergedPairs) {
    $filledRef = $pair[0]
    $blankRef  = $pair[1]

    $filledCell = $ws.Range($filledRef)
    $blankCell  = $ws.Range($blankRef)
    $mergedRange = $ws.Range($filledRef + ":" + $blankRef)

    $value = $filledCell.Value()

    # Split the merged range back into individual cells.
    $mergedRange.UnMerge()

    # Give the newly-separated blank cell the same formatting as its
    # former merge partner, then fill in the matching value.
    $filledCell.Copy()
    $blankCell.PasteSpecial(-4122)
    $blankCell.Value = $value
}

$excel.CutCopyMode = $false

$ws.Range("G4").Select() | Out-Null
